# Seguimiento et3.xlsx - revisión de tareas 7 (fila 16) y 9 (fila 18):
# se añade el tiempo empleado (columna G) que faltaba, lo que ajusta los
# totales de la hoja y reduce ligeramente la altura automática de esas filas.
# También se deja la hoja situada/seleccionada donde se estaba trabajando.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# Tarea 7 (fila 16): se registra el tiempo empleado (Recurso asociado = 5)
$ws.Range("G16").Value = 5
$ws.Rows.Item(16).RowHeight = 13.8

# Tarea 9 (fila 18): se registra el tiempo empleado (Recurso asociado = 5)
$ws.Range("G18").Value = 5
$ws.Rows.Item(18).RowHeight = 13.8

# Ligero reajuste del ancho de la columna L (tiempo empleado total)
$ws.Columns.Item(12).ColumnWidth = 26.3

# Ligero reajuste del ancho de la columna A de la hoja de estados
$ws2.Columns.Item(1).ColumnWidth = 19.8

# Dejar la vista/selección donde estaba el usuario trabajando
$ws.Activate()
$ws.Range("A13").Select()
$ws.Range("B22").Select()
